# TS7-Vytvorenie rezervacie.xlsx
# "pridavam vyslednu technicku dokumentaciu" -- the empty "Dodatocne
# informacie:" column (old column E, a header with no actual row content)
# is removed from the test-scenario table; the former column F
# ("Ocakavany vysledok:", the expected-result text) shifts left to become
# the new column E, and the CEVA logo picture that was anchored above the
# old column F now needs to sit above the new column E. The sheet is also
# switched to landscape for printing, and the selection is moved off the
# (now removed) old E15 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the logo picture's position (in points, relative to column F)
# before the column shift so it can be re-anchored afterwards -- deleting
# a column does not automatically drag floating pictures along with it.
$logo = $ws.Shapes.Item(1)
$oldAnchorCol = $ws.Range("F1")
$logoOffsetInCol = $logo.Left - $oldAnchorCol.Left
$logoTop = $logo.Top

# Drop the empty "Dodatočné informácie:" column; everything to its right
# (the "Očakávaný výsledok:" column and its header/body styling) shifts
# one column to the left.
$ws.Range("E:E").Delete()

# Re-anchor the logo above the new column E (the former column F), keeping
# the same in-column offset, top position, and size it always had.
$newAnchorCol = $ws.Range("E1")
$logo.Left = $newAnchorCol.Left + $logoOffsetInCol
$logo.Top = $logoTop

# Print the sheet in landscape now that it is narrower.
$ws.PageSetup.Orientation = 2   # xlLandscape

# Move the active selection (old selection pointed at the now-removed E15).
$ws.Range("J5").Select() | Out-Null
